$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Settings": clear the old "Cred_NA_Behorighetshantering_Berit"
# row (row 7, cols A:B) and move the selection to A7:B8.
# ---------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("A7").ClearContents() | Out-Null
$wsSettings.Range("B7").ClearContents() | Out-Null
$wsSettings.Range("A7:B8").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Assets": drop the obsolete "URLBerit" row (old row 6), which
# shifts every following row up by one, then append the two new rows
# needed for the receipt-saving feature.
# ---------------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Activate() | Out-Null

$wsAssets.Rows.Item(6).Delete() | Out-Null

$wsAssets.Range("A6").Value = "OePHämtaStatus"
$wsAssets.Range("B6").Value = "OePHämtaStatusArvoderadeUppdrag"

# Duplicate row 13's formatting (style "s=2") down onto the new row 14,
# then overwrite its contents.
$wsAssets.Rows.Item(13).Copy() | Out-Null
$wsAssets.Rows.Item(14).Insert() | Out-Null
$wsAssets.Rows.Item(14).RowHeight = 14.25
$wsAssets.Range("A14").Value = "FolderPathToCopies"
$wsAssets.Range("B14").Value = "FolderPathToCopies"

# Row 15 simply takes the plain (unstyled) defaults already on the sheet.
$wsAssets.Range("A15").Value = "FolderPathToReceipt"
$wsAssets.Range("B15").Value = "FolderPathToReceipt"

$wsAssets.Range("B15").Select() | Out-Null

Write-Output "done"
